# Update odds values in row 2 of Sheet1 to match the 2024-11-13 FlashScore refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "G2"  = 2.1
    "H2"  = 3.25
    "I2"  = 3.6
    "J2"  = 2.88
    "K2"  = 2.1
    "L2"  = 4
    "M2"  = 1.06
    "N2"  = 10
    "W2"  = 7
    "X2"  = 9.5
    "Y2"  = 9
    "Z2"  = 19
    "AA2" = 19
    "AD2" = 6.5
    "AG2" = 10
    "AH2" = 17
    "AI2" = 13
    "AJ2" = 41
    "AK2" = 29
    "AL2" = 41
    "AM2" = 4
    "AN2" = 12
    "AS2" = 8.5
    "AU2" = 5.5
    "AV2" = 21
    "AW2" = 29
    "AX2" = 67
    "AY2" = 101
    "BA2" = 301
    "BB2" = 151
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
